$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.331.28'
$ws.Range("E2").Value = '  +7.39%  '

$ws.Range("D3").Value = '2.673.77'
$ws.Range("E3").Value = '  +11.24%  '

$ws.Range("E4").Value = '  +0.48%  '

$ws.Range("D5").Value = '314.60'
$ws.Range("E5").Value = '  +7.40%  '

$ws.Range("D6").Value = '104.91'
$ws.Range("E6").Value = '  +12.89%  '

$ws.Range("D7").Value = '0.614'
$ws.Range("E7").Value = '  +10.50%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("E9").Value = '  +21.16%  '

$ws.Range("D10").Value = '41.27'
$ws.Range("E10").Value = '  +21.49%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '55.86'
$ws.Range("E11").Value = '  +4.70%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '0.0868'
$ws.Range("E12").Value = '  +12.25%  '

$ws.Range("D13").Value = '8.40'
$ws.Range("E13").Value = '  +20.72%  '

$ws.Range("D14").Value = '3.078.80'
$ws.Range("E14").Value = '  +11.25%  '

$ws.Range("E15").Value = '  +3.68%  '

$ws.Range("D16").Value = '2.697.88'
$ws.Range("E16").Value = '  +12.72%  '

$ws.Range("D17").Value = '0.949'
$ws.Range("E17").Value = '  +14.90%  '

$ws.Range("D18").Value = '15.49'
$ws.Range("E18").Value = '  +9.71%  '

$ws.Range("D19").Value = '48.535.72'
$ws.Range("E19").Value = '  +7.88%  '

$ws.Range("D20").Value = '0.0000104'
$ws.Range("E20").Value = '  +11.17%  '

$ws.Range("D21").Value = '13.33'
$ws.Range("E21").Value = '  +7.94%  '

$ws.Range("D22").Value = '6.91'
$ws.Range("E22").Value = '  +13.71%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '74.41'
$ws.Range("E23").Value = '  +12.05%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '285.04'
$ws.Range("E24").Value = '  +19.85%  '

$ws.Range("D25").Value = '3.12'
$ws.Range("E25").Value = '  +13.46%  '

$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  +18.41%  '

$ws.Range("D27").Value = '29.97'
$ws.Range("E27").Value = '  +42.84%  '

$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("E29").Value = '  +2.76%  '

$ws.Range("D30").Value = '10.78'
$ws.Range("E30").Value = '  +13.58%  '

$ws.Range("D31").Value = '40.72'
$ws.Range("E31").Value = '  +9.10%  '

$ws.Range("D32").Value = '2.33'
$ws.Range("E32").Value = '  +5.16%  '

$ws.Range("D33").Value = '6.24'
$ws.Range("E33").Value = '  +16.21%  '

$ws.Range("D34").Value = '3.73'
$ws.Range("E34").Value = '  -1.39%  '

$ws.Range("D35").Value = '0.0861'
$ws.Range("E35").Value = '  +14.03%  '

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.89'
$ws.Range("E36").Value = '  +6.97%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '2.25'
$ws.Range("E37").Value = '  +13.52%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '153.80'
$ws.Range("E38").Value = '  +3.75%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '0.126'
$ws.Range("E39").Value = '  +13.87%  '

$ws.Range("E40").Value = '  +9.92%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '4.36'
$ws.Range("E41").Value = '  +17.88%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '16.40'
$ws.Range("E42").Value = '  +14.76%  '

$ws.Range("D43").Value = '3.76'
$ws.Range("E43").Value = '  +20.28%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0341'
$ws.Range("E44").Value = '  +17.31%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '22.19'
$ws.Range("E45").Value = '  +42.17%  '

$ws.Range("D46").Value = '2.242.53'
$ws.Range("E46").Value = '  +13.55%  '

$ws.Range("D47").Value = '97.43'
$ws.Range("E47").Value = '  +10.19%  '

$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.25%  '

$ws.Range("D49").Value = '9.99'
$ws.Range("E49").Value = '  +19.66%  '

$ws.Range("D50").Value = '1.89'
$ws.Range("E50").Value = '  +11.25%  '

$ws.Range("D51").Value = '113.82'
$ws.Range("E51").Value = '  +14.64%  '
